$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Commentaires d'ordre génér" + (bookmark) + "al " -> single run
#    "Commentaires d'ordre général " and drop the old _GoBack bookmark that
#    used to sit between the two runs.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Commentaires d’ordre général ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Commentaires d’ordre général ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Merge the three runs "Projet_Gestion_Ecole_" / "FOUILLARD_NGUYEN" / ".zip"
#    into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Projet_Gestion_Ecole_FOUILLARD_NGUYEN.zip", $true, $false, $false, $false,
    $false, $true, 1, $false, "Projet_Gestion_Ecole_FOUILLARD_NGUYEN.zip", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Add Nam's first name after "Nguyen" as its own run, then place a
#    _GoBack bookmark right after it (mirrors Word parking _GoBack at the
#    most recent edit).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Nguyen") | Out-Null
$rng.Collapse(0)
# Insert a sentinel "X" after " Nam" so the bookmark insertion point is not
# exactly the paragraph-end boundary (inserting/bookmarking right at a
# paragraph end is unreliable), then remove the sentinel afterwards.
$rng.InsertAfter(" NamX")

# Give " Nam" (everything except the sentinel) its own run by toggling a
# character property on and back off, which keeps it from being re-merged
# into the previous "Nguyen" run.
$namOnly = $d.Range($rng.Start, $rng.End - 1)
$namOnly.Bold = 1
$namOnly.Bold = 0

# Place the _GoBack bookmark right before the sentinel (i.e. right after
# " Nam").
$xPos = $rng.End - 1
$bmRange = $d.Range($xPos, $xPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the sentinel character.
$sentinel = $d.Range($xPos, $xPos + 1)
$sentinel.Delete()
